$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 ("展览"): update a handful of "want to go" counts (column F),
# then insert two new rows (new events on 2024-09-15) before the existing
# "肥西·星域动漫游戏嘉年华" row, pushing it from row 9 down to row 11.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(2,6).Value = 5478
$ws1.Cells.Item(3,6).Value = 378
$ws1.Cells.Item(4,6).Value = 634
$ws1.Cells.Item(6,6).Value = 827
$ws1.Cells.Item(7,6).Value = 27
$ws1.Cells.Item(8,6).Value = 355

# Insert two blank rows at position 9 (existing row 9 and below shift down).
$ws1.Rows.Item(9).Resize(2).Insert()

# Copy the formatting of column A from row 8 onto the two new rows so the
# index cells keep the same bold/border/centered style (s="1").
$ws1.Cells.Item(8,1).Copy()
$ws1.Range("A9:A10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 9: 合肥·书香璃樱动漫游戏嘉年华
$ws1.Cells.Item(9,1).Value = 8
$ws1.Cells.Item(9,2).Value = "'2024-09-15"
$ws1.Cells.Item(9,3).Value = "合肥·书香璃樱动漫游戏嘉年华"
$ws1.Cells.Item(9,4).Value = "阜阳北路与金海路交口 格律诗宴会大酒店(北城店)"
$ws1.Cells.Item(9,5).Value = "2024.09.15 10:00-09.15 17:00"
$ws1.Cells.Item(9,6).Value = 0
$ws1.Cells.Item(9,7).Value = 50
$ws1.Cells.Item(9,8).Value = "https://show.bilibili.com/platform/detail.html?id=90735"
$ws1.Cells.Item(9,9).Value = "//i2.hdslb.com/bfs/openplatform/202408/7alsu0yg1723110506313.jpeg"

# Row 10: 合肥·曙光次元动漫游戏嘉年华
$ws1.Cells.Item(10,1).Value = 9
$ws1.Cells.Item(10,2).Value = "'2024-09-15"
$ws1.Cells.Item(10,3).Value = "合肥·曙光次元动漫游戏嘉年华"
$ws1.Cells.Item(10,4).Value = "田埠西路199号 吉祥如意宴会楼蜀山店"
$ws1.Cells.Item(10,5).Value = "2024.09.15 10:00-09.15 17:00"
$ws1.Cells.Item(10,6).Value = 0
$ws1.Cells.Item(10,7).Value = 50
$ws1.Cells.Item(10,8).Value = "https://show.bilibili.com/platform/detail.html?id=90733"
$ws1.Cells.Item(10,9).Value = "//i1.hdslb.com/bfs/openplatform/202408/bNZ6vKL01723113544322.jpeg"

# Row 11 (formerly row 9): 肥西·星域动漫游戏嘉年华 - only the index changes.
$ws1.Cells.Item(11,1).Value = 10

# ---------------------------------------------------------------------------
# Sheet 4 ("全部类型"): same count updates, then insert the same two new
# rows before the existing "肥西·星域动漫游戏嘉年华" row (was row 10, becomes
# row 12).
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Cells.Item(2,6).Value = 5478
$ws4.Cells.Item(3,6).Value = 378
$ws4.Cells.Item(4,6).Value = 634
$ws4.Cells.Item(6,6).Value = 827
$ws4.Cells.Item(7,6).Value = 27
$ws4.Cells.Item(9,6).Value = 355

# Insert two blank rows at position 10.
$ws4.Rows.Item(10).Resize(2).Insert()

$ws4.Cells.Item(9,1).Copy()
$ws4.Range("A10:A11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 10: 合肥·书香璃樱动漫游戏嘉年华
$ws4.Cells.Item(10,1).Value = 9
$ws4.Cells.Item(10,2).Value = "'2024-09-15"
$ws4.Cells.Item(10,3).Value = "合肥·书香璃樱动漫游戏嘉年华"
$ws4.Cells.Item(10,4).Value = "阜阳北路与金海路交口 格律诗宴会大酒店(北城店)"
$ws4.Cells.Item(10,5).Value = "2024.09.15 10:00-09.15 17:00"
$ws4.Cells.Item(10,6).Value = 0
$ws4.Cells.Item(10,7).Value = 50
$ws4.Cells.Item(10,8).Value = "https://show.bilibili.com/platform/detail.html?id=90735"
$ws4.Cells.Item(10,9).Value = "//i2.hdslb.com/bfs/openplatform/202408/7alsu0yg1723110506313.jpeg"

# Row 11: 合肥·曙光次元动漫游戏嘉年华
$ws4.Cells.Item(11,1).Value = 10
$ws4.Cells.Item(11,2).Value = "'2024-09-15"
$ws4.Cells.Item(11,3).Value = "合肥·曙光次元动漫游戏嘉年华"
$ws4.Cells.Item(11,4).Value = "田埠西路199号 吉祥如意宴会楼蜀山店"
$ws4.Cells.Item(11,5).Value = "2024.09.15 10:00-09.15 17:00"
$ws4.Cells.Item(11,6).Value = 0
$ws4.Cells.Item(11,7).Value = 50
$ws4.Cells.Item(11,8).Value = "https://show.bilibili.com/platform/detail.html?id=90733"
$ws4.Cells.Item(11,9).Value = "//i1.hdslb.com/bfs/openplatform/202408/bNZ6vKL01723113544322.jpeg"

# Row 12 (formerly row 10): 肥西·星域动漫游戏嘉年华 - only the index changes.
$ws4.Cells.Item(12,1).Value = 11

# Row 13 (formerly row 11): 合肥·《四月是你的谎言》... - only the index changes.
$ws4.Cells.Item(13,1).Value = 12

# Row 14 (formerly row 12): 合肥·一生必听的钢琴曲... - only the index changes.
$ws4.Cells.Item(14,1).Value = 13
